# Refresh the legacy GSC export data on the "Chart" sheet:
#   - the oldest day in the rolling window (2025-09-27, row 2) has aged out
#     of the export and is dropped, shifting every later date up one row
#   - "Not indexed" / "Indexed" crawl stats are not yet available for the
#     two most-recent days in the refreshed window, so those cells are
#     cleared while Impressions keeps its per-day history intact

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the row for 2025-09-27 -- everything below shifts up one row,
# shrinking the table from A1:D89 to A1:D88
$ws.Rows.Item(2).Delete()

# The two newest dates (now rows 2 and 3, i.e. 2025-09-28 and 2025-09-29)
# don't have "Not indexed" / "Indexed" coverage numbers yet in this refresh
$ws.Range("B2:C3").ClearContents()
